$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "G2" = 4.953047333333333
    "H2" = 14.859142
    "I2" = 0.7703204220313993
    "J2" = 0.7703204220313993
    "M2" = 1.918906333333333
    "N2" = 5.756718999999999
    "O2" = 0.006524019162508824
    "P2" = 0.006524019162508824
    "Q2" = 9.50443389723311
    "R2" = 85.53990507509798
    "S2" = 0.005025585194604733
    "T2" = 0.005025585194604733
    "G3" = 4.953047333333333
    "H3" = 14.859142
    "I3" = 0.7703204220313993
    "J3" = 0.7703204220313993
    "O3" = 0.6163557430885885
    "P3" = 0.6163557430885885
    "Q3" = 897.9299832578565
    "R3" = 8081.369849320708
    "S3" = 0.4747914161374782
    "T3" = 0.4747914161374782
    "G4" = 4.953047333333333
    "H4" = 14.859142
    "I4" = 0.7703204220313993
    "J4" = 0.7703204220313993
    "M4" = 29.04767233333333
    "N4" = 87.143017
    "O4" = 0.09875811426384234
    "P4" = 0.09875811426384236
    "Q4" = 143.8744959901571
    "R4" = 1294.870463911414
    "S4" = 0.07607539225874818
    "T4" = 0.0760753922587482
    "G5" = 4.953047333333333
    "H5" = 14.859142
    "I5" = 0.7703204220313993
    "J5" = 0.7703204220313993
    "M5" = 81.87450533333333
    "N5" = 245.623516
    "O5" = 0.2783621234850603
    "P5" = 0.2783621234850603
    "Q5" = 405.5283003092524
    "R5" = 3649.754702783272
    "S5" = 0.2144280284405681
    "T5" = 0.2144280284405681
    "I6" = 0.135969508894967
    "J6" = 0.135969508894967
    "M6" = 1.918906333333333
    "N6" = 5.756718999999999
    "O6" = 0.006524019162508824
    "P6" = 0.006524019162508824
    "Q6" = 1.677630726605333
    "R6" = 15.098676539448
    "S6" = 0.0008870676815476787
    "T6" = 0.0008870676815476787
    "I7" = 0.135969508894967
    "J7" = 0.135969508894967
    "O7" = 0.6163557430885885
    "P7" = 0.6163557430885885
    "S7" = 0.08380558769234783
    "T7" = 0.08380558769234783
    "I8" = 0.135969508894967
    "J8" = 0.135969508894967
    "M8" = 29.04767233333333
    "N8" = 87.143017
    "O8" = 0.09875811426384234
    "P8" = 0.09875811426384236
    "Q8" = 25.39533420482934
    "R8" = 228.558007843464
    "S8" = 0.01342809229584768
    "T8" = 0.01342809229584768
    "I9" = 0.135969508894967
    "J9" = 0.135969508894967
    "M9" = 81.87450533333333
    "N9" = 245.623516
    "O9" = 0.2783621234850603
    "P9" = 0.2783621234850603
    "Q9" = 71.57993253074133
    "R9" = 644.219392776672
    "S9" = 0.03784876122522381
    "T9" = 0.03784876122522381
    "G10" = 0.5382536666666667
    "H10" = 1.614761
    "I10" = 0.08371165542397027
    "J10" = 0.08371165542397027
    "M10" = 1.918906333333333
    "N10" = 5.756718999999999
    "O10" = 0.006524019162508824
    "P10" = 0.006524019162508824
    "Q10" = 1.032858369906556
    "R10" = 9.295725329159
    "S10" = 0.0005461364441113177
    "T10" = 0.0005461364441113177
    "G11" = 0.5382536666666667
    "H11" = 1.614761
    "I11" = 0.08371165542397027
    "J11" = 0.08371165542397027
    "O11" = 0.6163557430885885
    "P11" = 0.6163557430885885
    "Q11" = 97.57914135926825
    "R11" = 878.2122722334141
    "S11" = 0.05159615958401707
    "T11" = 0.05159615958401707
    "G12" = 0.5382536666666667
    "H12" = 1.614761
    "I12" = 0.08371165542397027
    "J12" = 0.08371165542397027
    "M12" = 29.04767233333333
    "N12" = 87.143017
    "O12" = 0.09875811426384234
    "P12" = 0.09875811426384236
    "Q12" = 15.63501614154856
    "R12" = 140.715145273937
    "S12" = 0.008267205231575854
    "T12" = 0.008267205231575855
    "G13" = 0.5382536666666667
    "H13" = 1.614761
    "I13" = 0.08371165542397027
    "J13" = 0.08371165542397027
    "M13" = 81.87450533333333
    "N13" = 245.623516
    "O13" = 0.2783621234850603
    "P13" = 0.2783621234850603
    "Q13" = 44.06925270218623
    "R13" = 396.623274319676
    "S13" = 0.02330215416426603
    "T13" = 0.02330215416426603
    "G14" = 0.06428833333333334
    "H14" = 0.192865
    "I14" = 0.009998413649663342
    "J14" = 0.009998413649663342
    "M14" = 1.918906333333333
    "N14" = 5.756718999999999
    "O14" = 0.006524019162508824
    "P14" = 0.006524019162508824
    "Q14" = 0.1233632899927778
    "R14" = 1.110269609935
    "S14" = 0.00006522984224509343
    "T14" = 0.00006522984224509343
    "G15" = 0.06428833333333334
    "H15" = 0.192865
    "I15" = 0.009998413649663342
    "J15" = 0.009998413649663342
    "O15" = 0.6163557430885885
    "P15" = 0.6163557430885885
    "Q15" = 11.65472853150111
    "R15" = 104.89255678351
    "S15" = 0.006162579674745336
    "T15" = 0.006162579674745336
    "G16" = 0.06428833333333334
    "H16" = 0.192865
    "I16" = 0.009998413649663342
    "J16" = 0.009998413649663342
    "M16" = 29.04767233333333
    "N16" = 87.143017
    "O16" = 0.09875811426384234
    "P16" = 0.09875811426384236
    "Q16" = 1.867426441522778
    "R16" = 16.806837973705
    "S16" = 0.0009874244776706133
    "T16" = 0.0009874244776706133
    "G17" = 0.06428833333333334
    "H17" = 0.192865
    "I17" = 0.009998413649663342
    "J17" = 0.009998413649663342
    "M17" = 81.87450533333333
    "N17" = 245.623516
    "O17" = 0.2783621234850603
    "P17" = 0.2783621234850603
    "Q17" = 5.263575490371111
    "R17" = 47.37217941334
    "S17" = 0.0027831796550023
    "T17" = 0.0027831796550023
}

foreach ($key in $values.Keys) {
    $ws.Range($key).Value = $values[$key]
}

Write-Output "done"